# Hata giderme: Yöne bağlı olarak hatalı hesaplama yapma sorunu giderildi.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C updates
$ws.Range("C2").Value = 582.29
$ws.Range("C3").Value = 487.9
$ws.Range("C4").Value = 393.51
$ws.Range("C5").Value = 299.13
$ws.Range("C6").Value = 204.75
$ws.Range("C7").Value = 110.37

$ws.Range("C14").Value = 769.12
$ws.Range("C15").Value = 358.6
$ws.Range("C16").Value = 59.82
$ws.Range("C17").Value = 59.88
$ws.Range("C18").Value = 59.95
$ws.Range("C19").Value = 60.02
$ws.Range("C20").Value = 60.08
$ws.Range("C21").Value = 60.15
$ws.Range("C22").Value = 60.22
$ws.Range("C23").Value = 60.28
$ws.Range("C24").Value = 141.87
$ws.Range("C25").Value = 284.71
$ws.Range("C27").Value = 715.54
$ws.Range("C28").Value = 109.01
$ws.Range("C29").Value = 55.48
$ws.Range("C30").Value = 55.7
$ws.Range("C31").Value = 55.93
$ws.Range("C32").Value = 56.16
$ws.Range("C33").Value = 56.38

# C34 changes from a numeric value to the "-" placeholder text
$ws.Range("C34").Value = "-"

# Column E updates
$ws.Range("E17").Value = 620.3
$ws.Range("E18").Value = 43.83
$ws.Range("E19").Value = -50.77
$ws.Range("E20").Value = -50.24
$ws.Range("E21").Value = -50.36
$ws.Range("E22").Value = -50.79
$ws.Range("E23").Value = -50.43
$ws.Range("E24").Value = -50.87
$ws.Range("E25").Value = -50.85
$ws.Range("E27").Value = -51.91
$ws.Range("E28").Value = -49.57
$ws.Range("E29").Value = -49.44
$ws.Range("E30").Value = -49.46
$ws.Range("E31").Value = -49.48
$ws.Range("E32").Value = -49.49
$ws.Range("E33").Value = -49.51
$ws.Range("E34").Value = -49.52
$ws.Range("E35").Value = -49.53
$ws.Range("E36").Value = -49.55
$ws.Range("E37").Value = -49.57

# E38 changes from the "-" placeholder text to a numeric value
$ws.Range("E38").Value = -49.23

$ws.Range("E39").Value = -49.84
$ws.Range("E40").Value = -50.45
$ws.Range("E41").Value = -201.72

# E42 changes from a numeric value to the "-" placeholder text
$ws.Range("E42").Value = "-"
